$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the R10 rule row.
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection on the sheet, matching the authored change.
$ws.Range("E8").Select()
